# Applies the cryptocurrency price/volume update described in the commit diff
# (GitHub Actions scheduled refresh of cryptos.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep Price/Volume (and Coin/Link) columns formatted as Text so that
# numeric-looking values (e.g. "318.06") stay strings, matching the
# original inlineStr cell type used throughout the sheet.
$ws.Range('B2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '45.443.45'
$ws.Range('E2').Value = '  +7.04%  '
$ws.Range('D3').Value = '2.378.85'
$ws.Range('E3').Value = '  +4.14%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '318.06'
$ws.Range('E5').Value = '  +2.34%  '
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '111.40'
$ws.Range('E6').Value = '  +7.69%  '
$ws.Range('E7').Value = '  +2.47%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +5.54%  '
$ws.Range('D10').Value = '42.04'
$ws.Range('E10').Value = '  +8.25%  '
$ws.Range('E11').Value = '  +3.43%  '
$ws.Range('D12').Value = '8.68'
$ws.Range('E12').Value = '  +5.49%  '
$ws.Range('E13').Value = '  +4.32%  '
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = '15.76'
$ws.Range('E15').Value = '  +4.87%  '
$ws.Range('D16').Value = '2.739.88'
$ws.Range('E16').Value = '  +4.00%  '
$ws.Range('D17').Value = '2.391.70'
$ws.Range('E17').Value = '  +4.70%  '
$ws.Range('D18').Value = '45.235.11'
$ws.Range('E18').Value = '  +6.70%  '
$ws.Range('D19').Value = '7.64'
$ws.Range('E19').Value = '  +5.36%  '
$ws.Range('E20').Value = '  +4.05%  '
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  -2.57%  '
$ws.Range('D22').Value = '75.17'
$ws.Range('E22').Value = '  +3.15%  '
$ws.Range('D23').Value = '3.56'
$ws.Range('E23').Value = '  +4.38%  '
$ws.Range('D24').Value = '269.22'
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  +7.96%  '
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('D27').Value = '7.60'
$ws.Range('E27').Value = '  +8.75%  '
$ws.Range('D28').Value = '11.30'
$ws.Range('E28').Value = '  +6.03%  '
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('D30').Value = '22.94'
$ws.Range('E30').Value = '  +3.34%  '
$ws.Range('D31').Value = '38.92'
$ws.Range('E31').Value = '  +9.47%  '
$ws.Range('E32').Value = '  +9.29%  '
$ws.Range('D33').Value = '169.84'
$ws.Range('E33').Value = '  +3.23%  '
$ws.Range('D34').Value = '3.01'
$ws.Range('E34').Value = '  +16.96%  '
$ws.Range('E35').Value = '  +2.54%  '
$ws.Range('E36').Value = '  +5.15%  '
$ws.Range('D37').Value = '4.84'
$ws.Range('E37').Value = '  +8.29%  '
$ws.Range('D38').Value = '3.10'
$ws.Range('E38').Value = '  +13.63%  '
$ws.Range('E39').Value = '  +5.33%  '
$ws.Range('E40').Value = '  +6.01%  '
$ws.Range('E41').Value = '  +11.41%  '
$ws.Range('D42').Value = '106.24'
$ws.Range('E42').Value = '  +7.11%  '
$ws.Range('D43').Value = '13.87'
$ws.Range('E43').Value = '  +16.53%  '
$ws.Range('E44').Value = '  +6.65%  '
$ws.Range('D45').Value = '71.78'
$ws.Range('E45').Value = '  +4.23%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').Value = '118.13'
$ws.Range('E47').Value = '  +7.46%  '
$ws.Range('D48').Value = '5.76'
$ws.Range('E48').Value = '  +11.60%  '
$ws.Range('E49').Value = '  +19.89%  '
$ws.Range('D50').Value = '79.51'
$ws.Range('E50').Value = '  +2.72%  '
$ws.Range('D51').Value = '9.16'
$ws.Range('E51').Value = '  +6.45%  '
